# Adds PASS/FAIL "Results" values to the TestCases and TestSteps sheets.
$wb = $excel.ActiveWorkbook

# --- TestCases sheet: column D ("Results") for rows 2-4 ---
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Range("D2").Value = "PASS"
$wsTestCases.Range("D3").Value = "PASS"
$wsTestCases.Range("D4").Value = "FAIL"

# --- TestSteps sheet: column H ("Results") for rows 2-33 ---
$wsTestSteps = $wb.Worksheets.Item("TestSteps")
for ($r = 2; $r -le 33; $r++) {
    $wsTestSteps.Cells.Item($r, 8).Value = "PASS"
}
